$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rId1 / sheet1.xml
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 997
$wsExpo.Range("F6").Value = 1158
$wsExpo.Range("F7").Value = 935
$wsExpo.Range("F8").Value = 284
$wsExpo.Range("F11").Value = 896
$wsExpo.Range("F12").Value = 322
$wsExpo.Range("F15").Value = 1376
$wsExpo.Range("F17").Value = 1280
$wsExpo.Range("F18").Value = 2940
$wsExpo.Range("F19").Value = 292
$wsExpo.Range("F20").Value = 1562
$wsExpo.Range("F21").Value = 1314
$wsExpo.Range("F22").Value = 757
$wsExpo.Range("F23").Value = 217
$wsExpo.Range("F24").Value = 1308
$wsExpo.Range("F26").Value = 1076
$wsExpo.Range("F28").Value = 3320
$wsExpo.Range("F31").Value = 1471

# Sheet "全部类型" (All types) - rId4 / sheet4.xml
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 997
$wsAll.Range("F10").Value = 1158
$wsAll.Range("F11").Value = 935
$wsAll.Range("F12").Value = 284
$wsAll.Range("F23").Value = 896
$wsAll.Range("F24").Value = 322
$wsAll.Range("F27").Value = 1376
$wsAll.Range("F29").Value = 1280
$wsAll.Range("F30").Value = 2940
$wsAll.Range("F31").Value = 292
$wsAll.Range("F32").Value = 1562
$wsAll.Range("F33").Value = 1314
$wsAll.Range("F34").Value = 757
$wsAll.Range("F35").Value = 217
$wsAll.Range("F36").Value = 1308
$wsAll.Range("F40").Value = 1076
$wsAll.Range("F42").Value = 3320
$wsAll.Range("F45").Value = 1471
